$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.069299666666667
$ws.Range("H2").Value = 3.207899
$ws.Range("I2").Value = 0.003616700200628781
$ws.Range("J2").Value = 0.003616700200628781
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 65.27176930503124
$ws.Range("R2").Value = 587.4459237452811
$ws.Range("S2").Value = 0.0007391137212867595
$ws.Range("T2").Value = 0.0007391137212867595
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.069299666666667
$ws.Range("H3").Value = 3.207899
$ws.Range("I3").Value = 0.003616700200628781
$ws.Range("J3").Value = 0.003616700200628781
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 113.6820230556447
$ws.Range("R3").Value = 1023.138207500802
$ws.Range("S3").Value = 0.001287293787171601
$ws.Range("T3").Value = 0.001287293787171601
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.069299666666667
$ws.Range("H4").Value = 3.207899
$ws.Range("I4").Value = 0.003616700200628781
$ws.Range("J4").Value = 0.003616700200628781
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 140.4401173206635
$ws.Range("R4").Value = 1263.961055885972
$ws.Range("S4").Value = 0.00159029269217042
$ws.Range("T4").Value = 0.00159029269217042
$ws.Range("I5").Value = 0.8238194745364892
$ws.Range("J5").Value = 0.8238194745364891
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 14867.73901845368
$ws.Range("R5").Value = 133809.6511660831
$ws.Range("S5").Value = 0.1683568567246209
$ws.Range("T5").Value = 0.1683568567246209
$ws.Range("I6").Value = 0.8238194745364892
$ws.Range("J6").Value = 0.8238194745364891
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.2932224493303102
$ws.Range("T6").Value = 0.2932224493303102
$ws.Range("I7").Value = 0.8238194745364892
$ws.Range("J7").Value = 0.8238194745364891
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 31989.74126604063
$ws.Range("R7").Value = 287907.6713943657
$ws.Range("S7").Value = 0.362240168481558
$ws.Range("T7").Value = 0.362240168481558
$ws.Range("G8").Value = 51.01955666666666
$ws.Range("H8").Value = 153.05867
$ws.Range("I8").Value = 0.1725638252628821
$ws.Range("J8").Value = 0.1725638252628821
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 3114.315693347859
$ws.Range("R8").Value = 28028.84124013073
$ws.Range("S8").Value = 0.03526537561154577
$ws.Range("T8").Value = 0.03526537561154577
$ws.Range("G9").Value = 51.01955666666666
$ws.Range("H9").Value = 153.05867
$ws.Range("I9").Value = 0.1725638252628821
$ws.Range("J9").Value = 0.1725638252628821
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 5424.116922573406
$ws.Range("R9").Value = 48817.05230316065
$ws.Range("S9").Value = 0.0614207227109545
$ws.Range("T9").Value = 0.06142072271095451
$ws.Range("G10").Value = 51.01955666666666
$ws.Range("H10").Value = 153.05867
$ws.Range("I10").Value = 0.1725638252628821
$ws.Range("J10").Value = 0.1725638252628821
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 6700.827417491861
$ws.Range("R10").Value = 60307.44675742675
$ws.Range("S10").Value = 0.07587772694038181
$ws.Range("T10").Value = 0.07587772694038181
